# Adding homework data and code.
#
# 1) The "datetimeFigureOut" date placeholders on the slide master and every
#    slide layout were re-cached from 8/24/2020 to 8/27/2020 (PowerPoint
#    recalculates these automatically whenever the deck is saved on a later
#    day).
# 2) On slide 2, the "Content Placeholder 2" text box had two adjacent runs
#    with identical formatting (". " and "(maybe a little R)") collapse into
#    a single run when the text box was re-edited.

$p = $ppt.ActivePresentation

# --- 1) Refresh the cached "today" date shown in the footer date fields ---

$newDate = "8/27/2020"

# Slide master's own Date Placeholder.
$masterShapes = $p.SlideMaster.Shapes
for ($i = 1; $i -le $masterShapes.Count; $i++) {
    $sh = $masterShapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

# Every slide layout's own Date Placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    $layoutShapes = $layout.Shapes
    for ($j = 1; $j -le $layoutShapes.Count; $j++) {
        $sh = $layoutShapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2) Merge ". " and "(maybe a little R)" into a single run on slide 2 ---

$slide2 = $p.Slides.Item(2)
$contentShape = $null
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $candidate = $slide2.Shapes.Item($i)
    if ($candidate.Name -eq "Content Placeholder 2") {
        $contentShape = $candidate
    }
}
if ($contentShape -eq $null) {
    $contentShape = $slide2.Shapes.Item(4)
}
$tr = $contentShape.TextFrame.TextRange

$fullText = $tr.Text
$fragment = ". (maybe a little R)"
$fragIndex = $fullText.IndexOf($fragment)
if ($fragIndex -ge 0) {
    $sub = $tr.Characters($fragIndex + 1, $fragment.Length)
    $sub.Text = $fragment
}
